# #7 fixed the slide dimension
#
# Widen the slide canvas (sldSz cx: 9144000 -> 12191969 EMU) and grow the
# single full-width content shape (TextBox on slide 2 / Picture on slides
# 3-16) on every content slide from cx="7315200" to cx="10058400" EMU so it
# still spans the slide. Heights / positions are untouched.

$p = $ppt.ActivePresentation

# --- 1. Presentation-level slide size -------------------------------------
# PowerPoint's object model expresses PageSetup sizes in points (1 pt =
# 12700 EMU). 12191969 EMU == 959.9975590551181 pt exactly.
$p.PageSetup.SlideWidth = 959.9975590551181

# --- 2. Per-slide content shape width --------------------------------------
# 10058400 EMU == 792 pt exactly.
$newWidth = 792

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Width -eq 576 -and $sh.Height -eq 360) {
            $sh.Width = $newWidth
        }
    }
}
